# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled update).
# Values are prefixed with a leading apostrophe so Excel stores them as literal
# text (matching the original inlineStr cells) instead of auto-converting
# numeric-looking strings (e.g. "113.19") into actual numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.256.21"
$ws.Range("E2").Value = "'  -1.22%  "
$ws.Range("D3").Value = "'2.268.63"
$ws.Range("E3").Value = "'  -1.57%  "
$ws.Range("E4").Value = "'  -0.15%  "
$ws.Range("D5").Value = "'113.19"
$ws.Range("E5").Value = "'  +4.08%  "
$ws.Range("D6").Value = "'264.51"
$ws.Range("E6").Value = "'  -2.64%  "
$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "'  -1.23%  "
$ws.Range("E8").Value = "'  +0.05%  "
$ws.Range("D9").Value = "'0.600"
$ws.Range("E9").Value = "'  -3.01%  "
$ws.Range("D10").Value = "'47.96"
$ws.Range("E10").Value = "'  +0.83%  "
$ws.Range("D11").Value = "'0.0924"
$ws.Range("E11").Value = "'  -1.66%  "
$ws.Range("D12").Value = "'8.76"
$ws.Range("E12").Value = "'  +3.48%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "'  -0.47%  "
$ws.Range("D14").Value = "'15.42"
$ws.Range("E14").Value = "'  -2.23%  "
$ws.Range("D15").Value = "'2.602.87"
$ws.Range("E15").Value = "'  -1.73%  "
$ws.Range("D16").Value = "'0.857"
$ws.Range("E16").Value = "'  -0.61%  "
$ws.Range("D17").Value = "'2.268.92"
$ws.Range("E17").Value = "'  -1.44%  "
$ws.Range("D18").Value = "'43.111.08"
$ws.Range("E18").Value = "'  -1.58%  "
$ws.Range("E19").Value = "'  -3.59%  "
$ws.Range("D20").Value = "'7.00"
$ws.Range("E20").Value = "'  +10.91%  "
$ws.Range("D21").Value = "'71.09"
$ws.Range("E21").Value = "'  -1.79%  "
$ws.Range("D22").Value = "'2.41"
$ws.Range("E22").Value = "'  -3.01%  "
$ws.Range("B23").Value = "'InternetComputer(DFINITY)"
$ws.Range("C23").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'9.84"
$ws.Range("E23").Value = "'  +6.22%  "
$ws.Range("B24").Value = "'BitcoinCash"
$ws.Range("C24").Value = "'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "'230.31"
$ws.Range("E24").Value = "'  -1.78%  "
$ws.Range("D25").Value = "'2.83"
$ws.Range("E25").Value = "'  -5.11%  "
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("D27").Value = "'11.32"
$ws.Range("E27").Value = "'  -0.64%  "
$ws.Range("D28").Value = "'3.87"
$ws.Range("E28").Value = "'  -1.74%  "
$ws.Range("D29").Value = "'41.12"
$ws.Range("E29").Value = "'  +0.29%  "
$ws.Range("E30").Value = "'  -2.27%  "
$ws.Range("E31").Value = "'  -1.36%  "
$ws.Range("D32").Value = "'171.56"
$ws.Range("E32").Value = "'  -3.46%  "
$ws.Range("D33").Value = "'21.30"
$ws.Range("E33").Value = "'  -2.86%  "
$ws.Range("D34").Value = "'0.0904"
$ws.Range("E34").Value = "'  -1.48%  "
$ws.Range("D35").Value = "'5.61"
$ws.Range("E35").Value = "'  -0.12%  "
$ws.Range("E36").Value = "'  -0.57%  "
$ws.Range("E37").Value = "'  -4.27%  "
$ws.Range("D38").Value = "'0.0350"
$ws.Range("E38").Value = "'  -2.15%  "
$ws.Range("D39").Value = "'3.78"
$ws.Range("E39").Value = "'  -0.69%  "
$ws.Range("D40").Value = "'0.104"
$ws.Range("E40").Value = "'  -8.89%  "
$ws.Range("D41").Value = "'14.10"
$ws.Range("E41").Value = "'  +15.47%  "
$ws.Range("D42").Value = "'74.84"
$ws.Range("E42").Value = "'  +11.29%  "
$ws.Range("D43").Value = "'2.42"
$ws.Range("E43").Value = "'  +3.36%  "
$ws.Range("D44").Value = "'0.234"
$ws.Range("E44").Value = "'  -1.53%  "
$ws.Range("D45").Value = "'6.11"
$ws.Range("E45").Value = "'  +10.34%  "
$ws.Range("E46").Value = "'  +0.02%  "
$ws.Range("D47").Value = "'1.37"
$ws.Range("E47").Value = "'  -1.88%  "
$ws.Range("D48").Value = "'8.58"
$ws.Range("E48").Value = "'  -2.55%  "
$ws.Range("D49").Value = "'0.0989"
$ws.Range("E49").Value = "'  -2.97%  "
$ws.Range("D50").Value = "'100.58"
$ws.Range("E50").Value = "'  +1.10%  "
$ws.Range("E51").Value = "'  +0.28%  "
